$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell A1: drop its bold/fill/border formatting (plain default style) ---
$ws.Range("A1").ClearFormats()

# --- New values for A2:A17 (replacing the old single DocEntry value in A2) ---
$values = @(70082,73300,73301,73303,73304,73306,73361,73362,73373,73375,73377,73378,73421,73436,73449,73459)
for ($i = 0; $i -lt $values.Length; $i++) {
  $ws.Cells.Item($i + 2, 1).Value = $values[$i]
}

# --- Style groups ---
# Group A (no bottom border): rows 2,3,6,8,10,13,14,15,16,17
$noBorderRows = @(2,3,6,8,10,13,14,15,16,17)
# Group B (with thin bottom border, same as the header's old look): rows 4,5,7,9,11,12
$borderRows = @(4,5,7,9,11,12)

$noBorderRange = $ws.Range("A2")
foreach ($r in $noBorderRows) {
  $noBorderRange = $excel.Union($noBorderRange, $ws.Cells.Item($r, 1))
}
$noBorderRange.ClearFormats()
$noBorderRange.Font.Bold = $true
$noBorderRange.Interior.Color = 65535

$borderRange = $ws.Cells.Item($borderRows[0], 1)
foreach ($r in $borderRows) {
  $borderRange = $excel.Union($borderRange, $ws.Cells.Item($r, 1))
}
$borderRange.ClearFormats()
$borderRange.Font.Bold = $true
$borderRange.Interior.Color = 65535
$borderRange.Borders.Item(9).LineStyle = 1
$borderRange.Borders.Item(9).Color = 4888374

# --- Selection / view state to match the target ---
$ws.Range("A2:A17").Select()
